$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target data for rows 2-10 (columns A-T).
# Column A: Sending cluster, B: Ligand symbol (Igf2), C: Receptor symbol (Igf1r),
# D: Target cluster, E-T: numeric stats. The "ECs" sending/target cluster is
# newly added (per Dr Hou's advice), expanding the 2x3 (FAPs/sCs) matrix to a
# full 3x3 (ECs/FAPs/sCs) matrix.
$data = @(
        @("ECs", "Igf2", "Igf1r", "ECs", 2, 0.6666666666666666, 13.24090133333333, 39.722704, 0.1214410874295642, 0.1214410874295642, 3, 1, 29.57110033333333, 88.713301, 0.5311523066901919, 0.5311523066901919, 391.5480218317671, 3523.932196485904, 0.0645037137151783, 0.0645037137151783),
        @("ECs", "Igf2", "Igf1r", "FAPs", 2, 0.6666666666666666, 13.24090133333333, 39.722704, 0.1214410874295642, 0.1214410874295642, 3, 1, 12.180337, 36.541011, 0.2187816490048282, 0.2187816490048282, 161.2786404237494, 1451.507763813744, 0.02656908136477958, 0.02656908136477958),
        @("ECs", "Igf2", "Igf1r", "sCs", 2, 0.6666666666666666, 13.24090133333333, 39.722704, 0.1214410874295642, 0.1214410874295642, 3, 1, 13.92204833333333, 41.76614499999999, 0.2500660443049799, 0.2500660443049799, 184.3404683395644, 1659.06421505608, 0.03036829234960635, 0.03036829234960635),
        @("FAPs", "Igf2", "Igf1r", "ECs", 3, 1, 5.270503666666666, 15.811511, 0.04833928449947708, 0.04833928449947708, 3, 1, 29.57110033333333, 88.713301, 0.5311523066901919, 0.5311523066901919, 155.8545927342012, 1402.691334607811, 0.02567552246565069, 0.02567552246565069),
        @("FAPs", "Igf2", "Igf1r", "FAPs", 3, 1, 5.270503666666666, 15.811511, 0.04833928449947708, 0.04833928449947708, 3, 1, 12.180337, 36.541011, 0.2187816490048282, 0.2187816490048282, 64.19651081973568, 577.7685973776211, 0.01057574837450913, 0.01057574837450913),
        @("FAPs", "Igf2", "Igf1r", "sCs", 3, 1, 5.270503666666666, 15.811511, 0.04833928449947708, 0.04833928449947708, 3, 1, 13.92204833333333, 41.76614499999999, 0.2500660443049799, 0.2500660443049799, 73.37620678834388, 660.3858610950949, 0.01208801365931726, 0.01208801365931726),
        @("sCs", "Igf2", "Igf1r", "ECs", 3, 1, 90.52007366666668, 271.560221, 0.8302196280709587, 0.8302196280709586, 3, 1, 29.57110033333333, 88.713301, 0.5311523066901919, 0.5311523066901919, 2676.778180577725, 24091.00362519952, 0.4409730705093629, 0.4409730705093629),
        @("sCs", "Igf2", "Igf1r", "FAPs", 3, 1, 90.52007366666668, 271.560221, 0.8302196280709587, 0.8302196280709586, 3, 1, 12.180337, 36.541011, 0.2187816490048282, 0.2187816490048282, 1102.565002524826, 9923.085022723433, 0.1816368192655395, 0.1816368192655395),
        @("sCs", "Igf2", "Igf1r", "sCs", 3, 1, 90.52007366666668, 271.560221, 0.8302196280709587, 0.8302196280709586, 3, 1, 13.92204833333333, 41.76614499999999, 0.2500660443049799, 0.2500660443049799, 1260.224840724227, 11342.02356651804, 0.2076097382960563, 0.2076097382960563)
)

$r = 2
foreach ($row in $data) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}
